# Auto-generated: update FFXIV leve-profit market data values
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 533.3333
$ws.Range("I2").Value = 450
$ws.Range("K2").Value = 450
$ws.Range("M2").Value = -337
$ws.Range("H39").Value = 345.4
$ws.Range("I39").Value = 240.83333
$ws.Range("J39").Value = 502.25
$ws.Range("K39").Value = 722.49999
$ws.Range("L39").Value = 1506.75
$ws.Range("M39").Value = -426.49999
$ws.Range("N39").Value = -2098.75
$ws.Range("H113").Value = 3401.8333
$ws.Range("I113").Value = 2905
$ws.Range("J113").Value = 3501.2
$ws.Range("K113").Value = 2905
$ws.Range("L113").Value = 3501.2
$ws.Range("M113").Value = 349
$ws.Range("N113").Value = -10009.2
$ws.Range("H116").Value = 14288636
$ws.Range("I116").Value = 18184920
$ws.Range("J116").Value = 2266.6667
$ws.Range("K116").Value = 18184920
$ws.Range("L116").Value = 2266.6667
$ws.Range("M116").Value = -18181478
$ws.Range("N116").Value = -9150.6667
$ws.Range("H132").Value = 2048.8333
$ws.Range("I132").Value = 1936.7885
$ws.Range("J132").Value = 2777.125
$ws.Range("K132").Value = 5810.3655
$ws.Range("L132").Value = 8331.375
$ws.Range("M132").Value = -3280.3655
$ws.Range("N132").Value = -13391.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26358.666
$ws.Range("J32").Value = 43722.223
$ws.Range("L32").Value = 43722.223
$ws.Range("N32").Value = -44296.223
$ws.Range("H45").Value = 1043.2
$ws.Range("I45").Value = 1053
$ws.Range("K45").Value = 1053
$ws.Range("M45").Value = -676
$ws.Range("H74").Value = 1351.2
$ws.Range("I74").Value = 1358.2273
$ws.Range("J74").Value = 1331.875
$ws.Range("K74").Value = 1358.2273
$ws.Range("L74").Value = 1331.875
$ws.Range("M74").Value = -484.2273
$ws.Range("N74").Value = -3079.875
$ws.Range("H77").Value = 1351.2
$ws.Range("I77").Value = 1358.2273
$ws.Range("J77").Value = 1331.875
$ws.Range("K77").Value = 6791.136500000001
$ws.Range("L77").Value = 6659.375
$ws.Range("M77").Value = -2423.136500000001
$ws.Range("N77").Value = -15395.375
$ws.Range("H132").Value = 409731.38
$ws.Range("I132").Value = 445707.56
$ws.Range("K132").Value = 1337122.68
$ws.Range("M132").Value = -1334592.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 881.75
$ws.Range("I29").Value = 881.75
$ws.Range("K29").Value = 881.75
$ws.Range("M29").Value = -592.75
$ws.Range("H107").Value = 2789.5625
$ws.Range("I107").Value = 2642.2
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 2642.2
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = -722.1999999999998
$ws.Range("N107").Value = -8840
$ws.Range("H134").Value = 837572.4399999999
$ws.Range("I134").Value = 1114062.5
$ws.Range("J134").Value = 8102.3335
$ws.Range("K134").Value = 3342187.5
$ws.Range("L134").Value = 24307.0005
$ws.Range("M134").Value = -3339652.5
$ws.Range("N134").Value = -29377.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1803.4706
$ws.Range("I16").Value = 1523.2222
$ws.Range("J16").Value = 2118.75
$ws.Range("K16").Value = 1523.2222
$ws.Range("L16").Value = 2118.75
$ws.Range("M16").Value = -1236.2222
$ws.Range("N16").Value = -2692.75
$ws.Range("H58").Value = 1350.1765
$ws.Range("I58").Value = 1458.4445
$ws.Range("J58").Value = 932.5714
$ws.Range("K58").Value = 1458.4445
$ws.Range("L58").Value = 932.5714
$ws.Range("M58").Value = -1255.4445
$ws.Range("N58").Value = -1338.5714
$ws.Range("H93").Value = 18740.334
$ws.Range("I93").Value = 18740.334
$ws.Range("K93").Value = 18740.334
$ws.Range("M93").Value = -16868.334
$ws.Range("H113").Value = 1803.4706
$ws.Range("I113").Value = 1523.2222
$ws.Range("J113").Value = 2118.75
$ws.Range("K113").Value = 1523.2222
$ws.Range("L113").Value = 2118.75
$ws.Range("M113").Value = 646.7778000000001
$ws.Range("N113").Value = -6458.75
$ws.Range("H118").Value = 22828
$ws.Range("J118").Value = 22828
$ws.Range("L118").Value = 22828
$ws.Range("N118").Value = -26142
$ws.Range("H122").Value = 1668.7273
$ws.Range("I122").Value = 1668.7273
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5006.1819
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2556.1819
$ws.Range("N122").Value = $null
$ws.Range("H132").Value = 1438.3541
$ws.Range("I132").Value = 1213.7
$ws.Range("J132").Value = 2561.625
$ws.Range("K132").Value = 3641.1
$ws.Range("L132").Value = 7684.875
$ws.Range("M132").Value = -1111.1
$ws.Range("N132").Value = -12744.875
$ws.Range("H134").Value = 1538.2307
$ws.Range("I134").Value = 1448.3158
$ws.Range("J134").Value = 1782.2858
$ws.Range("K134").Value = 4344.9474
$ws.Range("L134").Value = 5346.857400000001
$ws.Range("M134").Value = -1809.9474
$ws.Range("N134").Value = -10416.8574
$ws.Range("H136").Value = 1350.1765
$ws.Range("I136").Value = 1458.4445
$ws.Range("J136").Value = 932.5714
$ws.Range("K136").Value = 4375.333500000001
$ws.Range("L136").Value = 2797.7142
$ws.Range("M136").Value = -1825.333500000001
$ws.Range("N136").Value = -7897.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 169.83333
$ws.Range("I18").Value = 139.81818
$ws.Range("K18").Value = 419.4545400000001
$ws.Range("M18").Value = -250.4545400000001
$ws.Range("H69").Value = 806
$ws.Range("I69").Value = 612
$ws.Range("J69").Value = 1000
$ws.Range("K69").Value = 1836
$ws.Range("L69").Value = 3000
$ws.Range("M69").Value = -1025
$ws.Range("N69").Value = -4622
$ws.Range("H72").Value = 806
$ws.Range("I72").Value = 612
$ws.Range("J72").Value = 1000
$ws.Range("K72").Value = 5508
$ws.Range("L72").Value = 9000
$ws.Range("M72").Value = -1452
$ws.Range("N72").Value = -17112
$ws.Range("H131").Value = 871.1900000000001
$ws.Range("I131").Value = 536
$ws.Range("J131").Value = 892.5851
$ws.Range("K131").Value = 1608
$ws.Range("L131").Value = 2677.7553
$ws.Range("M131").Value = 3432
$ws.Range("N131").Value = -12757.7553
$ws.Range("H132").Value = 1817.1515
$ws.Range("J132").Value = 1975.3636
$ws.Range("L132").Value = 17778.2724
$ws.Range("N132").Value = -22838.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 29990
$ws.Range("J6").Value = 29990
$ws.Range("L6").Value = 29990
$ws.Range("N6").Value = -30216
$ws.Range("H16").Value = 29990
$ws.Range("J16").Value = 29990
$ws.Range("L16").Value = 29990
$ws.Range("N16").Value = -30490
$ws.Range("H122").Value = 4169.3945
$ws.Range("I122").Value = 3942.8708
$ws.Range("J122").Value = 5172.5713
$ws.Range("K122").Value = 11828.6124
$ws.Range("L122").Value = 15517.7139
$ws.Range("M122").Value = -9378.6124
$ws.Range("N122").Value = -20417.7139
$ws.Range("H126").Value = 4213.2
$ws.Range("I126").Value = 2799.3333
$ws.Range("J126").Value = 4566.6665
$ws.Range("K126").Value = 8397.999899999999
$ws.Range("L126").Value = 13699.9995
$ws.Range("M126").Value = -5927.999899999999
$ws.Range("N126").Value = -18639.9995
$ws.Range("H132").Value = 1909.25
$ws.Range("I132").Value = 1241.7826
$ws.Range("J132").Value = 4979.6
$ws.Range("K132").Value = 3725.3478
$ws.Range("L132").Value = 14938.8
$ws.Range("M132").Value = -1195.3478
$ws.Range("N132").Value = -19998.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 5000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = $null
$ws.Range("N19").Value = -5340
$ws.Range("H122").Value = 6989458.5
$ws.Range("I122").Value = 8221724.5
$ws.Range("K122").Value = 24665173.5
$ws.Range("M122").Value = -24662723.5
$ws.Range("H132").Value = 4711.316
$ws.Range("I132").Value = 4500.7334
$ws.Range("J132").Value = 5501
$ws.Range("K132").Value = 13502.2002
$ws.Range("L132").Value = 16503
$ws.Range("M132").Value = -10972.2002
$ws.Range("N132").Value = -21563
$ws.Range("H136").Value = 4216.351
$ws.Range("I136").Value = 4160.3335
$ws.Range("K136").Value = 12481.0005
$ws.Range("M136").Value = -9931.000499999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872
$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360
$ws.Range("H132").Value = 1828.4103
$ws.Range("I132").Value = 1477.7188
$ws.Range("K132").Value = 4433.1564
$ws.Range("M132").Value = -1903.1564
$ws.Range("H136").Value = 1742.878
$ws.Range("I136").Value = 1818.8
$ws.Range("K136").Value = 5456.4
$ws.Range("M136").Value = -2906.4

